# DQA FY2018 Q3 pivot table: clear out the sample/test facility rows and
# figures that were only there for demonstration purposes, leaving the
# header row and the blank pivot-table skeleton (with its original
# per-cell formatting) ready for real data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: clear the values in the data rows but keep whatever per-cell
# formatting (s=...) each cell already carries (this preserves the
# special formatting used on column K, on A5:B5, B7, C9:C10, G10:H10, etc.)
$ws.Range("A2:N11").ClearContents()

# Step 2: for the remaining cells that must go back to completely plain/
# unformatted + empty (i.e. they should not be written out as explicit
# cell records at all), clear both content and formatting.
$ws.Range("A2:J2").Clear()
$ws.Range("L2:N2").Clear()

$ws.Range("A3:N3").Clear()

$ws.Range("A4:J4").Clear()
$ws.Range("L4:N4").Clear()

$ws.Range("C5:J5").Clear()
$ws.Range("L5:N5").Clear()

$ws.Range("A6:J6").Clear()
$ws.Range("L6:N6").Clear()

$ws.Range("A7").Clear()
$ws.Range("C7:J7").Clear()
$ws.Range("L7:N7").Clear()

$ws.Range("A8:J8").Clear()
$ws.Range("L8:N8").Clear()

$ws.Range("A9:B9").Clear()
$ws.Range("D9:J9").Clear()
$ws.Range("L9:N9").Clear()

$ws.Range("A10:B10").Clear()
$ws.Range("D10:F10").Clear()
$ws.Range("I10:J10").Clear()
$ws.Range("L10:N10").Clear()

$ws.Range("A11:E11").Clear()
$ws.Range("L11:N11").Clear()

# Move the active selection back up to D1 (it had drifted to D16).
$null = $ws.Range("D1").Select()
